# Auto-generated script applying the Phantom_Profits value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 41.53846
$ws.Range("I6").Value = 44.090908
$ws.Range("J6").Value = 27.5
$ws.Range("K6").Value = 132.272724
$ws.Range("L6").Value = 82.5
$ws.Range("M6").Value = -20.27272399999998
$ws.Range("N6").Value = -306.5
$ws.Range("H17").Value = 1855.1666
$ws.Range("J17").Value = 1855.1666
$ws.Range("L17").Value = 5565.4998
$ws.Range("N17").Value = -5901.4998
$ws.Range("H32").Value = 3619.9092
$ws.Range("J32").Value = 4128.1665
$ws.Range("L32").Value = 4128.1665
$ws.Range("N32").Value = -4780.1665
$ws.Range("H42").Value = 189.75
$ws.Range("I42").Value = 159.5
$ws.Range("K42").Value = 478.5
$ws.Range("M42").Value = -248.5
$ws.Range("H43").Value = 5690.143
$ws.Range("J43").Value = 5915
$ws.Range("L43").Value = 5915
$ws.Range("N43").Value = -6053
$ws.Range("H58").Value = 299.2857
$ws.Range("I58").Value = 182.5
$ws.Range("K58").Value = 547.5
$ws.Range("M58").Value = -397.5
$ws.Range("H111").Value = 2176.6875
$ws.Range("J111").Value = 2434
$ws.Range("L111").Value = 7302
$ws.Range("N111").Value = -13436
$ws.Range("H121").Value = 996
$ws.Range("J121").Value = 996
$ws.Range("L121").Value = 2988
$ws.Range("N121").Value = -6482
$ws.Range("H132").Value = 4046.625
$ws.Range("I132").Value = 4179.4443
$ws.Range("K132").Value = 12538.3329
$ws.Range("M132").Value = -10008.3329
$ws.Range("H138").Value = 3003.4
$ws.Range("J138").Value = 3859.04
$ws.Range("L138").Value = 11577.12
$ws.Range("N138").Value = -21857.12
$ws.Range("H141").Value = 9007.941000000001
$ws.Range("I141").Value = 8953.916999999999
$ws.Range("K141").Value = 26861.751
$ws.Range("M141").Value = -21681.751

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2230.3684
$ws.Range("I45").Value = 2075.7222
$ws.Range("J45").Value = 5014
$ws.Range("K45").Value = 2075.7222
$ws.Range("L45").Value = 5014
$ws.Range("M45").Value = -1698.7222
$ws.Range("N45").Value = -5768
$ws.Range("H76").Value = 4762.6665
$ws.Range("J76").Value = 4762.6665
$ws.Range("L76").Value = 4762.6665
$ws.Range("N76").Value = -5438.6665
$ws.Range("H79").Value = 4762.6665
$ws.Range("J79").Value = 4762.6665
$ws.Range("L79").Value = 4762.6665
$ws.Range("N79").Value = -7102.6665
$ws.Range("H131").Value = 61666.332
$ws.Range("J131").Value = 61666.332
$ws.Range("L131").Value = 61666.332
$ws.Range("N131").Value = -71746.33199999999
$ws.Range("H132").Value = 4684.8125
$ws.Range("I132").Value = 4738.7095
$ws.Range("K132").Value = 14216.1285
$ws.Range("M132").Value = -11686.1285

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 80500
$ws.Range("J115").Value = 80500
$ws.Range("L115").Value = 80500
$ws.Range("N115").Value = -83634
$ws.Range("H134").Value = 3073.7
$ws.Range("I134").Value = 2569.2222
$ws.Range("K134").Value = 7707.6666
$ws.Range("M134").Value = -5172.6666

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 25542
$ws.Range("I50").Value = 25542
$ws.Range("K50").Value = 25542
$ws.Range("M50").Value = -24917
$ws.Range("H56").Value = 41083.332
$ws.Range("I56").Value = 10125
$ws.Range("K56").Value = 10125
$ws.Range("M56").Value = -9280
$ws.Range("H99").Value = 1679.6666
$ws.Range("I99").Value = 1679.6666
$ws.Range("K99").Value = 1679.6666
$ws.Range("M99").Value = -181.6666
$ws.Range("H122").Value = 3191.9092
$ws.Range("I122").Value = 3321.1
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 9963.299999999999
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -7513.299999999999
$ws.Range("N122").Value = -10600
$ws.Range("H126").Value = 1679.6666
$ws.Range("I126").Value = 1679.6666
$ws.Range("K126").Value = 5038.9998
$ws.Range("M126").Value = -2568.9998
$ws.Range("H134").Value = 7520899.5
$ws.Range("I134").Value = 8405270
$ws.Range("J134").Value = 3749.5
$ws.Range("K134").Value = 25215810
$ws.Range("L134").Value = 11248.5
$ws.Range("M134").Value = -25213275
$ws.Range("N134").Value = -16318.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1174262.1
$ws.Range("I4").Value = 646825.6
$ws.Range("K4").Value = 1940476.8
$ws.Range("M4").Value = -1940364.8
$ws.Range("H5").Value = 4577.087
$ws.Range("I5").Value = 3135
$ws.Range("K5").Value = 9405
$ws.Range("M5").Value = -9293
$ws.Range("H46").Value = 1500.5
$ws.Range("I46").Value = 833
$ws.Range("K46").Value = 2499
$ws.Range("M46").Value = -2408
$ws.Range("H75").Value = 2335.3333
$ws.Range("I75").Value = 2249
$ws.Range("K75").Value = 6747
$ws.Range("M75").Value = -5749
$ws.Range("H78").Value = 2335.3333
$ws.Range("I78").Value = 2249
$ws.Range("K78").Value = 20241
$ws.Range("M78").Value = -15249
$ws.Range("H114").Value = 1750
$ws.Range("I114").Value = 1000
$ws.Range("J114").Value = 2500
$ws.Range("K114").Value = 3000
$ws.Range("L114").Value = 7500
$ws.Range("M114").Value = 254
$ws.Range("N114").Value = -14008
$ws.Range("H117").Value = 2489.4285
$ws.Range("I117").Value = 2499.3333
$ws.Range("K117").Value = 7497.999899999999
$ws.Range("M117").Value = -4055.999899999999
$ws.Range("H121").Value = 90909900
$ws.Range("I121").Value = 859.75
$ws.Range("J121").Value = 142857940
$ws.Range("K121").Value = 2579.25
$ws.Range("L121").Value = 428573820
$ws.Range("M121").Value = -1269.25
$ws.Range("N121").Value = -428576440
$ws.Range("H129").Value = 2986.5
$ws.Range("J129").Value = 2986.5
$ws.Range("L129").Value = 8959.5
$ws.Range("N129").Value = -18959.5
$ws.Range("H135").Value = 4577.087
$ws.Range("I135").Value = 3135
$ws.Range("K135").Value = 28215
$ws.Range("M135").Value = -25680

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 398.93332
$ws.Range("J97").Value = 384.16666
$ws.Range("L97").Value = 384.16666
$ws.Range("N97").Value = -1376.16666
$ws.Range("H113").Value = 3670.3333
$ws.Range("I113").Value = 3505.5
$ws.Range("K113").Value = 3505.5
$ws.Range("M113").Value = -1335.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5999
$ws.Range("J7").Value = 5999
$ws.Range("L7").Value = 5999
$ws.Range("N7").Value = -6223
$ws.Range("H22").Value = 1236.375
$ws.Range("I22").Value = 1248.8
$ws.Range("K22").Value = 1248.8
$ws.Range("M22").Value = -953.8
$ws.Range("H27").Value = 1236.375
$ws.Range("I27").Value = 1248.8
$ws.Range("K27").Value = 1248.8
$ws.Range("M27").Value = -1141.8
$ws.Range("H55").Value = 1582.375
$ws.Range("I55").Value = 700
$ws.Range("K55").Value = 700
$ws.Range("M55").Value = -527
$ws.Range("H68").Value = 6250
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 10500
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 10500
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -11998
$ws.Range("H71").Value = 6250
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 10500
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 52500
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -59988
$ws.Range("H120").Value = 50698
$ws.Range("J120").Value = 50698
$ws.Range("L120").Value = 50698
$ws.Range("N120").Value = -60374
$ws.Range("H122").Value = 2227.8572
$ws.Range("I122").Value = 2227.8572
$ws.Range("K122").Value = 6683.571599999999
$ws.Range("M122").Value = -4233.571599999999
$ws.Range("H126").Value = 5999
$ws.Range("J126").Value = 5999
$ws.Range("L126").Value = 17997
$ws.Range("N126").Value = -22937
$ws.Range("H132").Value = 3750.7856
$ws.Range("I132").Value = 3868
$ws.Range("J132").Value = 3662.875
$ws.Range("K132").Value = 11604
$ws.Range("L132").Value = 10988.625
$ws.Range("M132").Value = -9074
$ws.Range("N132").Value = -16048.625

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 55000
$ws.Range("I29").Value = 37500
$ws.Range("K29").Value = 37500
$ws.Range("M29").Value = -37210
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -21108
$ws.Range("H107").Value = 2717.5806
$ws.Range("I107").Value = 2911.3333
$ws.Range("J107").Value = 2449.3076
$ws.Range("K107").Value = 8733.999899999999
$ws.Range("L107").Value = 7347.9228
$ws.Range("M107").Value = -6813.999899999999
$ws.Range("N107").Value = -11187.9228
$ws.Range("H122").Value = 2938.6296
$ws.Range("I122").Value = 3001.6191
$ws.Range("J122").Value = 2718.1667
$ws.Range("K122").Value = 9004.8573
$ws.Range("L122").Value = 8154.500100000001
$ws.Range("M122").Value = -6554.8573
$ws.Range("N122").Value = -13054.5001
$ws.Range("H136").Value = 3575.1924
$ws.Range("I136").Value = 3602
$ws.Range("K136").Value = 10806
$ws.Range("M136").Value = -8256
